# Workbook reference
$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- Rename the "Requested quantity" headers ---------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet right after "Monthly Trend" ----
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# --- Header row -----------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold / border / alignment) from an existing header
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$wsForecast.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats -- reuse same bold header formats

# --- Data rows --------------------------------------------------------
$data = @(
    @(45494.99999999999, 100, -40.57215237583957, 248.8152212891668),
    @(45557.99999999999, 167, 33.97106621270604, 303.35487437468),
    @(45585.99999999999, 196, 67.12767581705707, 328.3238722859797),
    @(45592.99999999999, 204, 68.04307711390807, 332.4720077916624),
    @(45599.99999999999, 211, 76.95433657455526, 340.9099215504094),
    @(45606.99999999999, 218, 78.12027271105346, 348.8949811040343),
    @(45613.99999999999, 226, 97.69760257096796, 350.89147028677),
    @(45620.99999999999, 233, 98.09355553612072, 360.7698169697106),
    @(45627.99999999999, 240, 110.4697558299891, 371.1351906270161),
    @(45634.99999999999, 248, 114.1441912550647, 370.7343732119126),
    @(45641.99999999999, 255, 114.2406591710208, 389.4437488175507),
    @(45648.99999999999, 263, 127.7504101570581, 394.2253876232255),
    @(45655.99999999999, 270, 144.6036252873211, 401.9747210131426)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Copy the date-format style from column A of Weekly Quantity down column A of the new sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats

Write-Host "done"
